{"js": "// Replace two-digit multiplication problems with new values, per commit diff.\nconst replacements = [\n  [\"22\u00d727=\", \"95\u00d779=\"],\n  [\"59\u00d728=\", \"26\u00d739=\"],\n  [\"62\u00d725=\", \"58\u00d764=\"],\n  [\"60\u00d728=\", \"28\u00d784=\"],\n  [\"64\u00d755=\", \"66\u00d752=\"],\n  [\"62\u00d731=\", \"80\u00d743=\"],\n  [\"81\u00d734=\", \"36\u00d763=\"],\n  [\"82\u00d783=\", \"27\u00d744=\"],\n  [\"52\u00d754=\", \"80\u00d720=\"],\n  [\"45\u00d724=\", \"95\u00d764=\"],\n  [\"46\u00d739=\", \"30\u00d751=\"],\n  [\"34\u00d751=\", \"82\u00d767=\"],\n  [\"63\u00d762=\", \"28\u00d767=\"],\n  [\"60\u00d780=\", \"57\u00d754=\"],\n  [\"48\u00d741=\", \"37\u00d737=\"],\n  [\"29\u00d765=\", \"18\u00d784=\"],\n  [\"69\u00d749=\", \"56\u00d750=\"],\n  [\"95\u00d721=\", \"27\u00d733=\"],\n  [\"36\u00d760=\", \"31\u00d733=\"],\n  [\"46\u00d741=\", \"50\u00d780=\"],\n  [\"80\u00d788=\", \"29\u00d762=\"],\n  [\"32\u00d793=\", \"25\u00d753=\"],\n  [\"31\u00d764=\", \"44\u00d732=\"],\n  [\"60\u00d716=\", \"75\u00d760=\"],\n  [\"42\u00d799=\", \"41\u00d746=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace two-digit multiplication problems with new values, per commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"22\u00d727=\", \"95\u00d779=\"),\n    @(\"59\u00d728=\", \"26\u00d739=\"),\n    @(\"62\u00d725=\", \"58\u00d764=\"),\n    @(\"60\u00d728=\", \"28\u00d784=\"),\n    @(\"64\u00d755=\", \"66\u00d752=\"),\n    @(\"62\u00d731=\", \"80\u00d743=\"),\n    @(\"81\u00d734=\", \"36\u00d763=\"),\n    @(\"82\u00d783=\", \"27\u00d744=\"),\n    @(\"52\u00d754=\", \"80\u00d720=\"),\n    @(\"45\u00d724=\", \"95\u00d764=\"),\n    @(\"46\u00d739=\", \"30\u00d751=\"),\n    @(\"34\u00d751=\", \"82\u00d767=\"),\n    @(\"63\u00d762=\", \"28\u00d767=\"),\n    @(\"60\u00d780=\", \"57\u00d754=\"),\n    @(\"48\u00d741=\", \"37\u00d737=\"),\n    @(\"29\u00d765=\", \"18\u00d784=\"),\n    @(\"69\u00d749=\", \"56\u00d750=\"),\n    @(\"95\u00d721=\", \"27\u00d733=\"),\n    @(\"36\u00d760=\", \"31\u00d733=\"),\n    @(\"46\u00d741=\", \"50\u00d780=\"),\n    @(\"80\u00d788=\", \"29\u00d762=\"),\n    @(\"32\u00d793=\", \"25\u00d753=\"),\n    @(\"31\u00d764=\", \"44\u00d732=\"),\n    @(\"60\u00d716=\", \"75\u00d760=\"),\n    @(\"42\u00d799=\", \"41\u00d746=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute(\n        $findText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $replaceText,\n        2\n    )\n}\n"}
